# "Made mode for months + days, got time to sync with PC on compile"
#
# 1. Row 2's feature changes from "Decimals on all time units" to
#    "Properly display time with 0s" (the old string becomes unused and
#    is dropped from the shared-string table; the new one is appended).
# 2. Rows 3-6 (Scroll through units with rotary / Days in months array /
#    Turn off/on screen / Update only when string size changes) are now
#    marked Done by putting "Yes" in column B.
# 3. The saved selection moves from I13 to A3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Yes"
$ws.Range("B4").Value = "Yes"
$ws.Range("B5").Value = "Yes"
$ws.Range("B6").Value = "Yes"

$ws.Range("A2").Value = "Properly display time with 0s"

[void]$ws.Range("A3").Select()
